$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old watch-list rows (7 through 16), collapsing the table
# back down to a single contiguous block (A1:B6).
$ws.Rows("7:16").Delete()

# The new tickers in rows 3-6 lose the banded "Normal_Sheet1" cell style
# that the old rows had (only row 2 keeps it), so clear formatting there
# before writing the new symbols.
$ws.Range("A3:A6").ClearFormats()

# New watch list symbols (note: row 4 is filled in after rows 2,3,5,6 so
# that the shared-string table order matches GME, MSFT, NVDA, AMZN, AAPL).
$ws.Range("A2").Value = "GME"
$ws.Range("A3").Value = "MSFT"
$ws.Range("A5").Value = "NVDA"
$ws.Range("A6").Value = "AMZN"
$ws.Range("A4").Value = "AAPL"

# Refresh the watch date for every remaining row.
$ws.Range("B2:B6").Value = 45475

# Move the active selection to A4, matching the saved view state.
[void]$ws.Range("A4").Select()
